$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert the year header row (E1:BL1) from text labels like "1960 [YR1960]"
# into plain numeric year values (1960-2019), left aligned.
$year = 1960
for ($col = 5; $col -le 64; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $year
    $cell.HorizontalAlignment = -4131
    $year++
}

# Update the selected cell / view
$ws.Range("G5").Select()
